# Daily automated update of the EPEX Spot / Gaz / CO2 price tracker.
# Adds the new day's column (Prix Spot) / rows (Gaz, CO2).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Prix Spot": append a new date column AW for 01-aug,
# carrying the same header formatting as the preceding AV column.
# ---------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Range("AV1").Copy($wsSpot.Range("AW1")) | Out-Null
$wsSpot.Range("AW1").Value = "01-aug"

$spotValues = @{
    2  = 102.22
    3  = 93.47
    4  = 81
    5  = 72.55
    6  = 63.98
    7  = 73.6
    8  = 83.14
    9  = 100.62
    10 = 100
    11 = 65.06
    12 = 45.81
    13 = 35.08
    14 = 31.87
    15 = 25.06
    16 = 20.01
    17 = 13.43
    18 = 17.42
    19 = 27.47
    20 = 44.8
    21 = 70.12
    22 = 87.34
    23 = 103.02
    24 = 102.52
    25 = 95.41
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Cells.Item($row, "AW").Value = $spotValues[$row]
}

# ---------------------------------------------------------------
# Sheet "Gaz": append 2025-07-30 as a new row.
# ---------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A46").NumberFormat = "@"
$wsGaz.Range("A46").Value = "2025-07-30"
$wsGaz.Range("A46").Style = "Normal"
$wsGaz.Range("B46").Value = 34.175

# ---------------------------------------------------------------
# Sheet "CO2": append 2025-07-30 as a new row.
# ---------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A46").NumberFormat = "@"
$wsCo2.Range("A46").Value = "2025-07-30"
$wsCo2.Range("A46").Style = "Normal"
$wsCo2.Range("B46").Value = 72.12
